$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell L2: 98.282 -> 98.3 ---
$ws.Range("L2").Value = 98.3

# --- New row 11 ---
$ws.Range("A10").Copy()
$ws.Range("A11").PasteSpecial(-4122)
$ws.Range("A11").Value = "2021年"
$ws.Range("B11").Value = 100.4
$ws.Range("C11").Value = 99.7
$ws.Range("D11").Value = 100.6
$ws.Range("F11").Value = 103.9
$ws.Range("G11").Value = 119.1
$ws.Range("H11").Value = 116.1
$ws.Range("I11").Value = 99.59999999999999
$ws.Range("J11").Value = 100.5
$ws.Range("K11").Value = 100.2
$ws.Range("L11").Value = 108.1
$ws.Range("M11").Value = 117.5
$ws.Range("N11").Value = 98.40000000000001
$ws.Range("O11").Value = 101.7
$ws.Range("P11").Value = 122.7
$ws.Range("Q11").Value = 113.1
$ws.Range("R11").Value = 101.4
$ws.Range("S11").Value = 103.2
$ws.Range("T11").Value = 101.1
$ws.Range("U11").Value = 99.59999999999999
$ws.Range("V11").Value = 100.6
$ws.Range("W11").Value = 145.1
$ws.Range("X11").Value = 105.1
$ws.Range("Y11").Value = 100.2
$ws.Range("Z11").Value = 104.4
$ws.Range("AA11").Value = 99.7
$ws.Range("AB11").Value = 128.2
$ws.Range("AC11").Value = 138.7
$ws.Range("AD11").Value = 104.1
$ws.Range("AE11").Value = 99.90000000000001
$ws.Range("AF11").Value = 99.90000000000001
$ws.Range("AG11").Value = 101.4
$ws.Range("AH11").Value = 104.9
$ws.Range("AI11").Value = 101.6
$ws.Range("AJ11").Value = 99.7
$ws.Range("AK11").Value = 106.7
$ws.Range("AL11").Value = 100.6
$ws.Range("AM11").Value = 103.7
$ws.Range("AN11").Value = 103.1
$ws.Range("AO11").Value = 101.8
$ws.Range("AP11").Value = 128.5
$ws.Range("AQ11").Value = 131
$ws.Range("E10").Copy()
$ws.Range("E11").PasteSpecial(-4122)

# --- New row 12 ---
$ws.Range("A10").Copy()
$ws.Range("A12").PasteSpecial(-4122)
$ws.Range("A12").Value = "2022年"
$ws.Range("B12").Value = 100.7413
$ws.Range("C12").Value = 101.5392
$ws.Range("D12").Value = 102.1574
$ws.Range("F12").Value = 104.7353
$ws.Range("G12").Value = 107.6503
$ws.Range("H12").Value = 104.1173
$ws.Range("I12").Value = 100.3221
$ws.Range("J12").Value = 100.9705
$ws.Range("K12").Value = 101.7138
$ws.Range("L12").Value = 104.1254
$ws.Range("M12").Value = 103.1563
$ws.Range("N12").Value = 98.9834
$ws.Range("O12").Value = 103.2351
$ws.Range("P12").Value = 105.4445
$ws.Range("Q12").Value = 108.3083
$ws.Range("R12").Value = 101.7589
$ws.Range("S12").Value = 101.2981
$ws.Range("T12").Value = 101.2704
$ws.Range("U12").Value = 100.1893
$ws.Range("V12").Value = 100.6374
$ws.Range("W12").Value = 117.0156
$ws.Range("X12").Value = 115.8872
$ws.Range("Y12").Value = 108.5743
$ws.Range("Z12").Value = 104.5709
$ws.Range("AA12").Value = 101.7804
$ws.Range("AB12").Value = 123.574
$ws.Range("AC12").Value = 135.896
$ws.Range("AD12").Value = 103.6331
$ws.Range("AE12").Value = 101.2861
$ws.Range("AF12").Value = 100.6655
$ws.Range("AG12").Value = 101.1452
$ws.Range("AH12").Value = 100.629
$ws.Range("AI12").Value = 100.9817
$ws.Range("AJ12").Value = 102.0971
$ws.Range("AK12").Value = 102.2142
$ws.Range("AL12").Value = 101.4779
$ws.Range("AM12").Value = 101.1858
$ws.Range("AN12").Value = 105.8546
$ws.Range("AO12").Value = 103.7321
$ws.Range("AP12").Value = 94.18989999999999
$ws.Range("AQ12").Value = 84.6237
$ws.Range("E10").Copy()
$ws.Range("E12").PasteSpecial(-4122)

$excel.CutCopyMode = $false
